$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Price column cells are plain text in the source data (e.g. "37.201.90" or
# "246.19"); forcing NumberFormat to text before assignment keeps Excel from
# reinterpreting numeric-looking strings as actual numbers, then the style is
# reset back to Normal so no extra formatting is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.201.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.002.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.848"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.294.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.993.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.158.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +14.16%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0662"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.64%  "
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.373.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.28%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.77%  "
